# Feature: duplicate "Planilha2" into a new "Planilha3" sheet (multi-row
# selection / print OP for PCP module), placed after Planilha2, and make it
# the active sheet.

$wb = $excel.ActiveWorkbook

# Source sheet to duplicate.
$srcSheet = $wb.Worksheets.Item("Planilha2")

# Copy it to the end of the workbook (after the current last sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)

# The copy becomes the new last sheet; rename it and make it active.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Planilha3"
$newSheet.Select()
$newSheet.Range("Q4").Select()
